# Apply the change described by the diff:
#  - Item "زولا رقبه" (row 32 of the "نواقص الأصناف" shortage report) was
#    removed from the source data. The report was regenerated, so every
#    row below it moves up one position (its item/price/etc. data shifts
#    up) while the "م" sequence numbers in column A stay put (1,2,3,...),
#    the grand total at the bottom is reduced by the removed row's price,
#    and the generated-on timestamp in the footer is refreshed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Deleting the whole worksheet row reproduces exactly what the source
# report generator does: every cell below shifts up one row (values,
# styles, merged cells, and the shared-string table all renumber
# consistently), including column A's sequence number.
$ws.Rows(32).Delete()

# The "م" column is a plain 1..N item counter tied to the row's position
# in the table, not to the underlying item - restore the original
# sequence (26..33) that the row shift above bumped up by one.
$ws.Cells.Item(32, 1).Value = 26
$ws.Cells.Item(33, 1).Value = 27
$ws.Cells.Item(34, 1).Value = 28
$ws.Cells.Item(35, 1).Value = 29
$ws.Cells.Item(36, 1).Value = 30
$ws.Cells.Item(37, 1).Value = 31
$ws.Cells.Item(38, 1).Value = 32
$ws.Cells.Item(39, 1).Value = 33

# The grand-total cell (now row 40 after the shift) is a hard-coded
# figure, not a formula, so it needs to drop by the removed row's price
# (50.00) by hand: 2193.32 -> 2143.32.
$ws.Cells.Item(40, 16).Value = 2143.3200000000002

# Refresh the "generated on" timestamp in the footer (now row 41).
$ws.Cells.Item(41, 1).Value = "Monday, 22 September, 2025 3:57 PM"
